$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.482.05'
$ws.Range('E2').Value = '  -3.65%  '
$ws.Range('D3').Value = '1.957.32'
$ws.Range('E3').Value = '  -2.31%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.011'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.30%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '321.76'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.009'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.37%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4765'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -4.81%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4064'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.68%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '53.28'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.41%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08444'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -6.47%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.060'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.08%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.15'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -4.89%  '
$ws.Range('D13').Value = '1.974.65'
$ws.Range('E13').Value = '  -6.17%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.634'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -5.20%  '
$ws.Range('E15').Value = '  -4.18%  '
$ws.Range('E16').Value = '  -0.28%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001074'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.67%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '89.28'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -5.41%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06629'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.75%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.74'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.58%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.009'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.40%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.821'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.70%  '
$ws.Range('D23').Value = '28.508.00'
$ws.Range('E23').Value = '  -3.75%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.61'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.289'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.75%  '
$ws.Range('D26').Value = '2.193.51'
$ws.Range('E26').Value = '  -4.25%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '154.23'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.83%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.24'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.32%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.004'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.38%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.166'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.73%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '123.93'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9892'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -6.35%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09608'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.46%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.450'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -7.23%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.604'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.99%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.661'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.69%  '
$ws.Range('E37').Value = '  -5.09%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06226'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.89%  '
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.800'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -5.20%  '
$ws.Range('E40').Value = '  -3.72%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6230'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.91%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.15'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.54%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1923'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.95%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.335'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.60%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5973'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.87%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '13.05'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.86%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.060'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.04%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.398'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.17%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00000000330'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.32%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06838'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.19%  '
